$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).Value = '63.538.46'
$ws.Cells.Item(2,5).Value = '  +0.54%  '

# Row 3
$ws.Cells.Item(3,4).Value = '2.636.86'
$ws.Cells.Item(3,5).Value = '  -0.94%  '

# Row 4
$ws.Cells.Item(4,4).Value = "'0.999"
$ws.Cells.Item(4,5).Value = '  -0.12%  '

# Row 5
$ws.Cells.Item(5,4).Value = "'606.99"
$ws.Cells.Item(5,5).Value = '  -0.11%  '

# Row 6
$ws.Cells.Item(6,4).Value = "'147.53"
$ws.Cells.Item(6,5).Value = '  +3.12%  '

# Row 7
$ws.Cells.Item(7,4).Value = "'0.999"
$ws.Cells.Item(7,5).Value = '  -0.14%  '

# Row 8
$ws.Cells.Item(8,5).Value = '  +0.20%  '

# Row 9
$ws.Cells.Item(9,4).Value = "'0.109"
$ws.Cells.Item(9,5).Value = '  +1.98%  '

# Row 10
$ws.Cells.Item(10,4).Value = "'0.379"
$ws.Cells.Item(10,5).Value = '  +6.46%  '

# Row 11
$ws.Cells.Item(11,4).Value = "'5.54"
$ws.Cells.Item(11,5).Value = '  -2.25%  '

# Row 12
$ws.Cells.Item(12,5).Value = '  -0.55%  '

# Row 13
$ws.Cells.Item(13,4).Value = "'27.39"
$ws.Cells.Item(13,5).Value = '  +0.42%  '

# Row 14
$ws.Cells.Item(14,4).Value = '3.103.58'
$ws.Cells.Item(14,5).Value = '  -1.25%  '

# Row 15
$ws.Cells.Item(15,4).Value = '63.366.88'
$ws.Cells.Item(15,5).Value = '  +0.46%  '

# Row 16
$ws.Cells.Item(16,4).Value = "'0.0000147"
$ws.Cells.Item(16,5).Value = '  +1.82%  '

# Row 17
$ws.Cells.Item(17,4).Value = '2.626.65'
$ws.Cells.Item(17,5).Value = '  -0.37%  '

# Row 18
$ws.Cells.Item(18,4).Value = "'11.67"
$ws.Cells.Item(18,5).Value = '  +1.92%  '

# Row 19
$ws.Cells.Item(19,4).Value = "'4.56"
$ws.Cells.Item(19,5).Value = '  +3.98%  '

# Row 20
$ws.Cells.Item(20,4).Value = "'344.93"
$ws.Cells.Item(20,5).Value = '  +1.66%  '

# Row 21
$ws.Cells.Item(21,4).Value = "'6.90"
$ws.Cells.Item(21,5).Value = '  +0.82%  '

# Row 22
$ws.Cells.Item(22,5).Value = '  -0.12%  '

# Row 23
$ws.Cells.Item(23,4).Value = "'5.63"
$ws.Cells.Item(23,5).Value = '  -2.47%  '

# Row 24
$ws.Cells.Item(24,4).Value = "'66.18"
$ws.Cells.Item(24,5).Value = '  -2.10%  '

# Row 25
$ws.Cells.Item(25,2).Value = 'SuiNetwork'
$ws.Cells.Item(25,3).Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Cells.Item(25,4).Value = "'1.63"
$ws.Cells.Item(25,5).Value = '  +4.56%  '

# Row 26
$ws.Cells.Item(26,2).Value = 'Fetch.AI'
$ws.Cells.Item(26,3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(26,4).Value = "'1.69"
$ws.Cells.Item(26,5).Value = '  +2.17%  '

# Row 27
$ws.Cells.Item(27,5).Value = '  +6.83%  '

# Row 28
$ws.Cells.Item(28,4).Value = "'556.16"
$ws.Cells.Item(28,5).Value = '  +3.12%  '

# Row 29
$ws.Cells.Item(29,2).Value = 'Aptos'
$ws.Cells.Item(29,3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(29,4).Value = "'8.07"
$ws.Cells.Item(29,5).Value = '  +2.86%  '

# Row 30
$ws.Cells.Item(30,2).Value = 'Kaspa'
$ws.Cells.Item(30,3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(30,4).Value = "'0.162"
$ws.Cells.Item(30,5).Value = '  -1.25%  '

# Row 31
$ws.Cells.Item(31,5).Value = '  -0.51%  '

# Row 32
$ws.Cells.Item(32,5).Value = '  +0.27%  '

# Row 33
$ws.Cells.Item(33,4).Value = '0.0₃0849'
$ws.Cells.Item(33,5).Value = '  +5.13%  '

# Row 34
$ws.Cells.Item(34,5).Value = '  -2.21%  '

# Row 35
$ws.Cells.Item(35,4).Value = "'5.34"
$ws.Cells.Item(35,5).Value = '  +5.07%  '

# Row 36
$ws.Cells.Item(36,4).Value = "'168.88"
$ws.Cells.Item(36,5).Value = '  -2.21%  '

# Row 37
$ws.Cells.Item(37,5).Value = '  +0.01%  '

# Row 38
$ws.Cells.Item(38,4).Value = "'0.404"
$ws.Cells.Item(38,5).Value = '  -0.24%  '

# Row 39
$ws.Cells.Item(39,4).Value = "'1.92"
$ws.Cells.Item(39,5).Value = '  +5.43%  '

# Row 40
$ws.Cells.Item(40,4).Value = "'19.07"
$ws.Cells.Item(40,5).Value = '  -0.75%  '

# Row 41
$ws.Cells.Item(41,5).Value = '  +0.02%  '

# Row 42
$ws.Cells.Item(42,4).Value = "'164.90"
$ws.Cells.Item(42,5).Value = '  -5.46%  '

# Row 43
$ws.Cells.Item(43,5).Value = '  -0.47%  '

# Row 44
$ws.Cells.Item(44,4).Value = "'3.78"
$ws.Cells.Item(44,5).Value = '  +1.17%  '

# Row 45
$ws.Cells.Item(45,4).Value = "'22.00"
$ws.Cells.Item(45,5).Value = '  -0.03%  '

# Row 46
$ws.Cells.Item(46,4).Value = "'0.0567"
$ws.Cells.Item(46,5).Value = '  +0.93%  '

# Row 47
$ws.Cells.Item(47,5).Value = '  -0.73%  '

# Row 48
$ws.Cells.Item(48,2).Value = 'dogwifhat'
$ws.Cells.Item(48,3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(48,4).Value = "'1.99"
$ws.Cells.Item(48,5).Value = '  +15.16%  '

# Row 49
$ws.Cells.Item(49,2).Value = 'VeChain'
$ws.Cells.Item(49,3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(49,4).Value = "'0.0244"
$ws.Cells.Item(49,5).Value = '  +2.23%  '

# Row 50
$ws.Cells.Item(50,4).Value = "'0.0955"
$ws.Cells.Item(50,5).Value = '  -0.70%  '

# Row 51
$ws.Cells.Item(51,4).Value = "'18.87"
$ws.Cells.Item(51,5).Value = '  +0.87%  '
